# Auto-generated: update cryptocurrency price/volume table to match
# the Fri Aug 23 13:35:17 UTC 2024 GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.839.74'
$ws.Range('E2').Value = '  +0.50%  '

$ws.Range('D3').Value = '2.637.16'
$ws.Range('E3').Value = '  +0.40%  '

$ws.Range('E4').Value = '  -0.01%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '581.91'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.03%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.79'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.10%  '

$ws.Range('E7').Value = '  +0.06%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.596'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.42%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '6.58'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.63%  '

$ws.Range('E10').Value = '  +2.03%  '

$ws.Range('E11').Value = '  +1.43%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.376'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.69%  '

$ws.Range('D13').Value = '3.112.32'
$ws.Range('E13').Value = '  +0.94%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.35'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +5.37%  '

$ws.Range('D15').Value = '60.810.02'
$ws.Range('E15').Value = '  +0.43%  '

$ws.Range('E16').Value = '  +1.51%  '

$ws.Range('D17').Value = '2.653.53'
$ws.Range('E17').Value = '  +1.27%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '11.60'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.58%  '

$ws.Range('E19').Value = '  +0.67%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '351.29'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.69%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.88'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.76%  '

$ws.Range('E22').Value = '  -0.07%  '

$ws.Range('E23').Value = '  +1.02%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '63.98'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.07%  '

$ws.Range('E25').Value = '  +1.85%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.998'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.18%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.38'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +5.52%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.00'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +7.53%  '

$ws.Range('D29').Value = '0.0₃0807'
$ws.Range('E29').Value = '  +1.15%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.76'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +5.92%  '

$ws.Range('B31').Value = 'Monero'
$ws.Range('C31').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '167.47'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.93%  '

$ws.Range('B32').Value = 'USDe'
$ws.Range('C32').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.998'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.04%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '19.94'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.23%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.60'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +8.16%  '

$ws.Range('E35').Value = '  +6.44%  '

$ws.Range('E36').Value = '  +6.86%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.68'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.66%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '341.48'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +8.44%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.10'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +5.24%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.903'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +7.12%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '38.31'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.00%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '138.08'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.18%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.29'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +4.86%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0573'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.42%  '

$ws.Range('B45').Value = 'InjectiveProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '21.02'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +3.87%  '

$ws.Range('B46').Value = 'Mantle'
$ws.Range('C46').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.622'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.35%  '

$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '20.26'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.73%  '

$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0250'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.57%  '

$ws.Range('E49').Value = '  +0.44%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.999'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.10%  '

$ws.Range('D51').Value = '2.087.47'
$ws.Range('E51').Value = '  +2.35%  '

